# Update ltcf case files - accommodate formatting of cases file for ltcfs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("routes")

# Correct region-code / region-name assignments for several zip rows
$ws.Range("C48").Value = "south_king_county"
$ws.Range("A50").Value = "East King County"
$ws.Range("A73").Value = "West King County"
$ws.Range("C73").Value = "west_king_county"
$ws.Range("C78").Value = "west_king_county"
$ws.Range("C79").Value = "west_king_county"
$ws.Range("C82").Value = "west_king_county"
$ws.Range("C83").Value = "west_king_county"
$ws.Range("A96").Value = "South King County"

# Move the active selection to E4, as in the saved workbook state
$ws.Range("E4").Select()
